$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F20").Value = -3.037240239421779
$ws.Range("F22").Value = -3.062305129253265
$ws.Range("F26").Value = -3.100133327727114
$ws.Range("F33").Value = -3.147488660590756
$ws.Range("F36").Value = -3.166768557658393
$ws.Range("F37").Value = -3.166986840509973
$ws.Range("F38").Value = -3.178121789665754
$ws.Range("F39").Value = -3.178727714492667
$ws.Range("F40").Value = -3.189086890847607
$ws.Range("F42").Value = -3.151612062211434
$ws.Range("F43").Value = -3.176933566071869
$ws.Range("F46").Value = -3.186905737946059
$ws.Range("F47").Value = -3.177952125946729
$ws.Range("F48").Value = -3.178212174314628
$ws.Range("F49").Value = -3.176700416238692
$ws.Range("F50").Value = -3.163038048458731
$ws.Range("F52").Value = -3.179434566913131
$ws.Range("F53").Value = -3.154234476745204
$ws.Range("F55").Value = -3.161839726021661
$ws.Range("F81").Value = -3.314300246440736
$ws.Range("F85").Value = -3.336685814913305
$ws.Range("F87").Value = -3.238388281621968
$ws.Range("F88").Value = -3.170137997571151
$ws.Range("F89").Value = -3.079341127588406
$ws.Range("F92").Value = -3.400957124031009
$ws.Range("F94").Value = -3.44475366772961
$ws.Range("F95").Value = -3.446981877464456
$ws.Range("F96").Value = -3.427035522429484
$ws.Range("F97").Value = -3.388681069861232
$ws.Range("F98").Value = -3.34719642309861
$ws.Range("F104").Value = -3.50886494137496
$ws.Range("F106").Value = -3.512105776710649
$ws.Range("F107").Value = -3.489116622696817
$ws.Range("F108").Value = -3.455432132214432
$ws.Range("F109").Value = -3.403021966738486
$ws.Range("F114").Value = -3.526000278900787
$ws.Range("F118").Value = -3.531048261035067
$ws.Range("F119").Value = -3.490412192081617
$ws.Range("F120").Value = -3.441504662442166
$ws.Range("F125").Value = -3.56816085410785
$ws.Range("F136").Value = -3.595040559483012
$ws.Range("F137").Value = -3.612148631224381
$ws.Range("F142").Value = -3.490481781365247
$ws.Range("F147").Value = -3.62059439912723
$ws.Range("F158").Value = -3.637767903514832
